$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Erick Silva" was renamed to "Erick da Silva" (Responsavel column, rows 3-14)
$ws.Range("B3:B14").Value = "Erick da Silva"

# Update the active selection to reflect the new working cell
$ws.Range("K16").Select()
